# Natmi following Dr Hou advice
# Update the LR-pairs sheet: rows 2-5 get refreshed values (with updated
# cluster labels) and six additional rows (6-11) are appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-D are text (cluster / gene symbol) labels; E-T are numeric values.
$rows = @(
  ,@("sCs","Nppc","Npr2","ECs", 1,0.3333333333333333,0.224617,0.673851,0.223433695150483,0.3014711375983185,2,1,17.947853,35.895706,0.4337138114960137,0.3586943164132274,4.031392897301,24.188357383806,0.09690627954035438,0.1081359836191469)
  ,@("sCs","Nppc","Npr2","FAPs",1,0.3333333333333333,0.224617,0.673851,0.223433695150483,0.3014711375983185,3,1,4.537805333333334,13.613416,0.1096570629894806,0.1360345147179691,1.019268220557333,9.173413985015999,0.02450108278308893,0.04101047990466135)
  ,@("sCs","Nppc","Npr2","M2",  1,0.3333333333333333,0.224617,0.673851,0.223433695150483,0.3014711375983185,1,0.3333333333333333,0.01289433333333333,0.038683,0.0003115943983216321,0.0003865468544291307,0.002896286470333333,0.026066578233,[double]"6.962068780519373E-05",0.0001165327199398017)
  ,@("sCs","Nppc","Npr2","sCs", 1,0.3333333333333333,0.224617,0.673851,0.223433695150483,0.3014711375983185,3,1,12.75897433333333,38.276923,0.3083234183436765,0.3824890567659189,2.865882537830333,25.792942840473,0.06888984066195586,0.1153094110621294)
  ,@("sCs","Nppc","Npr2","Nppc",1,0.3333333333333333,0.224617,0.673851,0.223433695150483,0.3014711375983185,2,1,6.124261000000001,12.248522,0.1479941127725076,0.1223955652484555,1.375613133037,8.253678798222001,0.03306687147727867,0.0368987302924411)
  ,@("Nppc","Nppc","Npr2","ECs", 2,1,0.780679,1.561358,0.7765663048495169,0.6985288624016814,2,1,17.947853,35.895706,0.4337138114960137,0.3586943164132274,14.011511932187,56.04604772874801,0.3368075319556593,0.2505583327940805)
  ,@("Nppc","Nppc","Npr2","FAPs",2,1,0.780679,1.561358,0.7765663048495169,0.6985288624016814,3,1,4.537805333333334,13.613416,0.1096570629894806,0.1360345147179691,3.542569329821334,21.255415978928,0.0851559802063917,0.09502403481330775)
  ,@("Nppc","Nppc","Npr2","M2",  2,1,0.780679,1.561358,0.7765663048495169,0.6985288624016814,1,0.3333333333333333,0.01289433333333333,0.038683,0.0003115943983216321,0.0003865468544291307,0.01006633525233333,0.060398011514,0.0002419737105164383,0.000270014134489329)
  ,@("Nppc","Nppc","Npr2","sCs", 2,1,0.780679,1.561358,0.7765663048495169,0.6985288624016814,3,1,12.75897433333333,38.276923,0.3083234183436765,0.3824890567659189,9.960663323572334,59.76397994143399,0.2394335776817206,0.2671796457037895)
  ,@("Nppc","Nppc","Npr2","Nppc",2,1,0.780679,1.561358,0.7765663048495169,0.6985288624016814,2,1,6.124261000000001,12.248522,0.1479941127725076,0.1223955652484555,4.781081953219001,19.124327812876,0.1149272412952289,0.0854968349560144)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $values[$c]
    }
}
